$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.614754438400269
$ws.Range("B1").Value = 2.771294355392456
$ws.Range("C1").Value = 3.226843118667603
$ws.Range("D1").Value = 3.60989236831665
$ws.Range("E1").Value = 1.538814663887024
